$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03273466666666667
$ws.Range("H2").Value = 0.098204
$ws.Range("I2").Value = 0.08359843399780884
$ws.Range("J2").Value = 0.08359843399780884
$ws.Range("M2").Value = 4.093680666666667
$ws.Range("N2").Value = 12.281042
$ws.Range("O2").Value = 0.1610908176055751
$ws.Range("P2").Value = 0.161090817605575
$ws.Range("Q2").Value = 0.1340052720631111
$ws.Range("R2").Value = 1.206047448568
$ws.Range("S2").Value = 0.01346694008325273
$ws.Range("T2").Value = 0.01346694008325273

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03273466666666667
$ws.Range("H3").Value = 0.098204
$ws.Range("I3").Value = 0.08359843399780884
$ws.Range("J3").Value = 0.08359843399780884
$ws.Range("O3").Value = 0.5606512265211691
$ws.Range("P3").Value = 0.5606512265211691
$ws.Range("Q3").Value = 0.466384249948
$ws.Range("R3").Value = 4.197458249532
$ws.Range("S3").Value = 0.04686956455612052
$ws.Range("T3").Value = 0.04686956455612052

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03273466666666667
$ws.Range("H4").Value = 0.098204
$ws.Range("I4").Value = 0.08359843399780884
$ws.Range("J4").Value = 0.08359843399780884
$ws.Range("M4").Value = 7.071161666666666
$ws.Range("N4").Value = 21.213485
$ws.Range("O4").Value = 0.2782579558732559
$ws.Range("P4").Value = 0.2782579558732559
$ws.Range("Q4").Value = 0.2314721201044444
$ws.Range("R4").Value = 2.08324908094
$ws.Range("S4").Value = 0.02326192935843558
$ws.Range("T4").Value = 0.02326192935843558

$ws.Range("I5").Value = 0.3399848984133119
$ws.Range("J5").Value = 0.3399848984133119
$ws.Range("M5").Value = 4.093680666666667
$ws.Range("N5").Value = 12.281042
$ws.Range("O5").Value = 0.1610908176055751
$ws.Range("P5").Value = 0.161090817605575
$ws.Range("Q5").Value = 0.544983519792
$ws.Range("R5").Value = 4.904851678128
$ws.Range("S5").Value = 0.05476844525894878
$ws.Range("T5").Value = 0.05476844525894878

$ws.Range("I6").Value = 0.3399848984133119
$ws.Range("J6").Value = 0.3399848984133119
$ws.Range("O6").Value = 0.5606512265211691
$ws.Range("P6").Value = 0.5606512265211691
$ws.Range("S6").Value = 0.1906129502940984
$ws.Range("T6").Value = 0.1906129502940984

$ws.Range("I7").Value = 0.3399848984133119
$ws.Range("J7").Value = 0.3399848984133119
$ws.Range("M7").Value = 7.071161666666666
$ws.Range("N7").Value = 21.213485
$ws.Range("O7").Value = 0.2782579558732559
$ws.Range("P7").Value = 0.2782579558732559
$ws.Range("Q7").Value = 0.9413696103599999
$ws.Range("R7").Value = 8.472326493240001
$ws.Range("S7").Value = 0.09460350286026471
$ws.Range("T7").Value = 0.09460350286026473

$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.2257076666666666
$ws.Range("H8").Value = 0.6771229999999999
$ws.Range("I8").Value = 0.5764166675888793
$ws.Range("J8").Value = 0.5764166675888793
$ws.Range("M8").Value = 4.093680666666667
$ws.Range("N8").Value = 12.281042
$ws.Range("O8").Value = 0.1610908176055751
$ws.Range("P8").Value = 0.161090817605575
$ws.Range("Q8").Value = 0.9239751113517777
$ws.Range("R8").Value = 8.315776002165999
$ws.Range("S8").Value = 0.09285543226337355
$ws.Range("T8").Value = 0.09285543226337353

$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.2257076666666666
$ws.Range("H9").Value = 0.6771229999999999
$ws.Range("I9").Value = 0.5764166675888793
$ws.Range("J9").Value = 0.5764166675888793
$ws.Range("O9").Value = 0.5606512265211691
$ws.Range("P9").Value = 0.5606512265211691
$ws.Range("Q9").Value = 3.215749892851
$ws.Range("R9").Value = 28.941749035659
$ws.Range("S9").Value = 0.3231687116709502
$ws.Range("T9").Value = 0.3231687116709502

$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.2257076666666666
$ws.Range("H10").Value = 0.6771229999999999
$ws.Range("I10").Value = 0.5764166675888793
$ws.Range("J10").Value = 0.5764166675888793
$ws.Range("M10").Value = 7.071161666666666
$ws.Range("N10").Value = 21.213485
$ws.Range("O10").Value = 0.2782579558732559
$ws.Range("P10").Value = 0.2782579558732559
$ws.Range("Q10").Value = 1.596015400406111
$ws.Range("R10").Value = 14.364138603655
$ws.Range("S10").Value = 0.1603925236545556
$ws.Range("T10").Value = 0.1603925236545556

